# Adds the raw and clean SSA data for June 23 (date serial 44005) across
# the relevant sheets of the workbook, mirroring the upstream commit
# "Raw and Clean Data from SSA for June 23th".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# out_vars: new row 24 (daily raw totals)
# ---------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")
$wsOut.Range("A23:J23").Copy($wsOut.Range("A24:J24"))
$wsOut.Range("A24").Value = 44005
$wsOut.Range("B24").Value = 191410
$wsOut.Range("C24").Value = 251355
$wsOut.Range("D24").Value = 59106
$wsOut.Range("E24").Value = 23377
$wsOut.Range("F24").Value = 31.598140118071154
$wsOut.Range("G24").Value = 60482
$wsOut.Range("H24").Value = 5219
$wsOut.Range("I24").Value = 5402
$wsOut.Range("J24").Value = 501871

# ---------------------------------------------------------------------
# dates_dx: new row 24 (daily observation counts by delay bucket)
# ---------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")
$wsDx.Range("A23:K23").Copy($wsDx.Range("A24:K24"))
$wsDx.Range("A24").Value = 44005
$wsDx.Range("A24").Style = $wsDx.Range("A22").Style
$wsDx.Range("B24").Value = 0
$wsDx.Range("C24").Value = 1
$wsDx.Range("D24").Value = 1
$wsDx.Range("E24").Value = 1
$wsDx.Range("F24").Value = 1
$wsDx.Range("G24").Value = 0
$wsDx.Range("H24").Value = 0
$wsDx.Range("I24").Value = 1
$wsDx.Range("J24").Value = 0
$wsDx.Range("K24").Value = 4

# ---------------------------------------------------------------------
# dates_sx: new row 24
# ---------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")
$wsSx.Range("A23:M23").Copy($wsSx.Range("A24:M24"))
$wsSx.Range("A24").Value = 44005
$wsSx.Range("B24").Value = 0
$wsSx.Range("C24").Value = 1
$wsSx.Range("D24").Value = 1
$wsSx.Range("E24").Value = 0
$wsSx.Range("F24").Value = 1
$wsSx.Range("G24").Value = 1
$wsSx.Range("H24").Value = 1
$wsSx.Range("I24").Value = 0
$wsSx.Range("J24").Value = 1
$wsSx.Range("K24").Value = 1
$wsSx.Range("L24").Value = 0
$wsSx.Range("M24").Value = 0

# ---------------------------------------------------------------------
# dates_deaths: new row 24
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("dates_deaths")
$wsDe.Range("A23:J23").Copy($wsDe.Range("A24:J24"))
$wsDe.Range("A24").Value = 44005
$wsDe.Range("B24").Value = 0
$wsDe.Range("C24").Value = 0
$wsDe.Range("D24").Value = 0
$wsDe.Range("E24").Value = 0
$wsDe.Range("F24").Value = 2
$wsDe.Range("G24").Value = 1
$wsDe.Range("H24").Value = 1
$wsDe.Range("I24").Value = 1
$wsDe.Range("J24").Value = 2

# ---------------------------------------------------------------------
# control_obs: new column X (44005) for the summary/control table
# ---------------------------------------------------------------------
$wsCo = $wb.Worksheets.Item("control_obs")

$wsCo.Range("W1").Copy($wsCo.Range("X1"))
$wsCo.Range("X1").Value = 44005

$wsCo.Range("X2").Value = 3753
$wsCo.Range("X3").Value = 3568
$wsCo.Range("X4").Value = 3568
$wsCo.Range("X5").Value = 3568
$wsCo.Range("X6").Value = 3568
$wsCo.Range("X7").Value = 2759
$wsCo.Range("X8").Value = 5390
$wsCo.Range("X10").Value = 163
$wsCo.Range("X11").Value = 163
$wsCo.Range("X12").Value = 163
$wsCo.Range("X13").Value = 163
$wsCo.Range("X14").Value = 163
$wsCo.Range("X15").Value = 98
$wsCo.Range("X16").Value = 175
$wsCo.Range("X18").Value = 870

$wsCo.Range("W20").Copy($wsCo.Range("X20"))
$wsCo.Range("X20").Formula = "=SUM(X2:X18)"

# ---------------------------------------------------------------------
# Restore a sensible selection on each sheet (mirrors the cursor moves
# made while editing); control_obs is reselected last so it remains the
# active tab, matching the workbook's saved view state.
# ---------------------------------------------------------------------
$wsOut.Activate()
$wsOut.Range("B26").Select()

$wsDx.Activate()
$wsDx.Range("K24").Select()

$wsSx.Activate()
$wsSx.Range("M24").Select()

$wsDe.Activate()
$wsDe.Range("K24").Select()

$wsCo.Activate()
$wsCo.Range("Z16").Select()
